
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Expand the "Main Control Signals" table with a new "Uses  Z?" column (BB)
# and convert the resSrc column (BA) from a single 0/1 bit into a 2-bit
# selector ("00"/"01"/"10") now that JAL needs a third result-source option.
# ---------------------------------------------------------------------------

# New "Uses  Z?" column values (BB27:BB33), written top-to-bottom first so
# the shared-string table picks up "N" / "Y" in the same order the original
# author entered them.
$ws.Range("AZ27").Copy() | Out-Null
$ws.Range("BB27").PasteSpecial(-4122) | Out-Null
$ws.Range("BB27").Value = "N"

$ws.Range("AZ28").Copy() | Out-Null
$ws.Range("BB28").PasteSpecial(-4122) | Out-Null
$ws.Range("BB28").Value = "N"

$ws.Range("AZ29").Copy() | Out-Null
$ws.Range("BB29").PasteSpecial(-4122) | Out-Null
$ws.Range("BB29").Value = "N"

$ws.Range("AZ30").Copy() | Out-Null
$ws.Range("BB30").PasteSpecial(-4122) | Out-Null
$ws.Range("BB30").Value = "N"

$ws.Range("AZ31").Copy() | Out-Null
$ws.Range("BB31").PasteSpecial(-4122) | Out-Null
$ws.Range("BB31").Value = "N"

$ws.Range("AZ32").Copy() | Out-Null
$ws.Range("BB32").PasteSpecial(-4122) | Out-Null
$ws.Range("BB32").Value = "Y"

$ws.Range("AZ33").Copy() | Out-Null
$ws.Range("BB33").PasteSpecial(-4122) | Out-Null
$ws.Range("BB33").Value = "N"

# BEQ's resSrc is "don't care" (the ALU result is discarded; only the Zero
# flag matters), recorded as text "dc".
$ws.Range("BA32").Value = "dc"

# New column header.
$ws.Range("BA26").Copy() | Out-Null
$ws.Range("BB26").PasteSpecial(-4122) | Out-Null
$ws.Range("BB26").Value = "Uses  Z?"

# Re-express resSrc (BA27:BA31) as an explicit 2-bit text code instead of a
# plain 0/1 bit, now that a third option ("10", PC+4 for JAL) exists.
$ws.Range("BA27").NumberFormat = "@"
$ws.Range("BA27").Value = "00"

$ws.Range("BA28").NumberFormat = "@"
$ws.Range("BA28").Value = "00"

$ws.Range("BA29").NumberFormat = "@"
$ws.Range("BA29").Value = "00"

$ws.Range("BA30").NumberFormat = "@"
$ws.Range("BA30").Value = "01"

$ws.Range("BA31").NumberFormat = "@"
$ws.Range("BA31").Value = "00"

# JAL (row 33) selects PC+4 as the result source -> resSrc = "10".
$ws.Range("BA33").Formula = "'10"

# ---------------------------------------------------------------------------
# Fill in the BEQ (row 32) and JAL (row 33) control-signal bits themselves.
# ---------------------------------------------------------------------------

# BEQ: pcSrc=1 (branch taken when ALU zero flag is set), regWr=0, aluSrc=0
# (compares two registers), memWr=0.
$ws.Range("AW32").Value = 1
$ws.Range("AX32").Value = 0
$ws.Range("AY32").Value = 0
$ws.Range("AZ32").Value = 0

# JAL: pcSrc=0, regWr=1 (writes return address to rd), aluSrc=1, memWr=0.
$ws.Range("AW33").Value = 0
$ws.Range("AX33").Value = 1
$ws.Range("AY33").Value = 1
$ws.Range("AZ33").Value = 0

# ---------------------------------------------------------------------------
# Note explaining BEQ's pcSrc behaviour.
# ---------------------------------------------------------------------------
$cmt = $ws.Range("AW32").AddComment()
$cmt.Text("Tom:" + [char]10 + "if zero is 1" + [char]10)

# ---------------------------------------------------------------------------
# Move the selection in the bottom-right (frozen) pane to reflect where the
# author ended up after making the edits above.
# ---------------------------------------------------------------------------
$ws.Range("M37").Select()
